$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the three Qualcomm FastConnect rows that were dropped this week ---
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(18).Delete()   # Qualcomm FastConnect ...1229 (Good Drivers table)
$ws.Rows.Item(16).Delete()   # Qualcomm FastConnect ...1277 (Good Drivers table)
$ws.Rows.Item(7).Delete()    # Qualcomm FastConnect ...1193 (Bad Drivers table)

# --- Updated figures for the "Bad Drivers" table ---
$ws.Range("C3").Value = 2222
$ws.Range("D3").Value = 87.6

$ws.Range("C4").Value = 103
$ws.Range("D4").Value = 95.6

$ws.Range("B5").Value = 14
$ws.Range("C5").Value = 2247
$ws.Range("D5").Value = 97

$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 300
$ws.Range("D6").Value = 98.2

# Totals row (shifted up to row 7 after the deletion above)
$ws.Range("B7").Value = 29
$ws.Range("C7").Value = 4872

# --- Updated figures for the "Good Drivers" table (rows shifted up after deletions) ---
$ws.Range("B15").Value = 449371   # Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4
$ws.Range("B18").Value = 331283   # Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8
$ws.Range("B25").Value = 77999    # Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9

# --- Narrow column A now that the Qualcomm entries (long names) are gone ---
$ws.Columns.Item(1).ColumnWidth = 45
